# Generate Report for Handback
# Updates the localization-status report to reflect a failed handback
# transform for the d07e63b8-... file in both the zh-cn and de-de
# localization sheets, and records the error detail message returned
# by the handback transform for each locale.

$wb = $excel.ActiveWorkbook

$statusFailed = "Handback transform failed"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusFailed
$wsOverview.Range("F3").Value = $statusFailed

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusFailed
$wsZhCn.Range("P3").Value = "Handback file name: ykthaay3.gwq is different with handoff file name: d07e63b8-7bf1-448a-b511-031f74b3232a.680998ba1e4137a2bb93ef21f8c4885c0e0cdb9f.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusFailed
$wsDeDe.Range("P3").Value = "Handback file name: ykthaay3.gwq is different with handoff file name: d07e63b8-7bf1-448a-b511-031f74b3232a.680998ba1e4137a2bb93ef21f8c4885c0e0cdb9f.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
